$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended to the "ScenarioResults" sheet (rows 519-533).
# Column A = scenario name, Column B = status, Column C = browser.
$newRows = @(
    @("Country with Parameter", "FAILED", "chrome"),
    @("Country with Parameter", "FAILED", "chrome"),
    @("Country with Parameter", "FAILED", "chrome"),
    @("Country with Parameter", "FAILED", "chrome"),
    @("Country with Parameter", "FAILED", "chrome"),
    @("Country with Parameter", "FAILED", "chrome"),
    @("Country with Parameter", "FAILED", "chrome"),
    @("Country with Parameter", "PASSED", "chrome"),
    @("Create a country", "PASSED", "chrome"),
    @("Create Nationality and Delete", "PASSED", "chrome"),
    @("Create a country", "PASSED", "chrome"),
    @("Create a country with base name and code (but generate unique)", "PASSED", "chrome"),
    @("Create Country", "PASSED", "chrome"),
    @("Create Nationality", "PASSED", "chrome"),
    @("Fee Functionality", "PASSED", "chrome")
)

$startRow = 519
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
